$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 308, shifting existing rows 308-315 down to 310-317
$ws.Rows("308:309").Insert()

# Fill new row 308
$ws.Range("A308").Value = 9
$ws.Range("B308").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C308").Value = "Metropolitana"
$ws.Range("D308").Value = 44448
$ws.Range("E308").Value = 13
$ws.Range("F308").Value = 100112040
$ws.Range("G308").Value = "Cilantro"
$ws.Range("H308").Value = "Sin especificar"
$ws.Range("I308").Value = "Primera"
$ws.Range("J308").Value = 52
$ws.Range("K308").Value = 6000
$ws.Range("L308").Value = 6000
$ws.Range("M308").Value = 6000
$ws.Range("N308").Value = "`$/caja 36 atados"
$ws.Range("O308").Value = "Región Metropolitana"
$ws.Range("P308").Value = 167
$ws.Range("Q308").Value = 36
$ws.Range("R308").Value = "Hortaliza"

# Fill new row 309
$ws.Range("A309").Value = 9
$ws.Range("B309").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C309").Value = "Metropolitana"
$ws.Range("D309").Value = 44448
$ws.Range("E309").Value = 13
$ws.Range("F309").Value = 100112040
$ws.Range("G309").Value = "Cilantro"
$ws.Range("H309").Value = "Sin especificar"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 133
$ws.Range("K309").Value = 10000
$ws.Range("L309").Value = 12000
$ws.Range("M309").Value = 11008
$ws.Range("N309").Value = "`$/docena de atados"
$ws.Range("O309").Value = "Región Metropolitana"
$ws.Range("P309").Value = 3669
$ws.Range("Q309").Value = 3
$ws.Range("R309").Value = "Hortaliza"
